$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# NumberFormat "@" + ClearFormats() keeps these as plain text cells
# (matching the original inline-string cells) without leaving a stray
# cell style behind, even for values that look numeric (e.g. "247.66").

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '36.878.95'
$c.ClearFormats()

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +4.28%  '
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.910.91'
$c.ClearFormats()

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +1.49%  '
$c.ClearFormats()

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.07%  '
$c.ClearFormats()

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '247.66'
$c.ClearFormats()

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +0.84%  '
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.685'
$c.ClearFormats()

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -0.89%  '
$c.ClearFormats()

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.ClearFormats()

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '47.45'
$c.ClearFormats()

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +9.75%  '
$c.ClearFormats()

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.375'
$c.ClearFormats()

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +5.81%  '
$c.ClearFormats()

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '58.15'
$c.ClearFormats()

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +6.08%  '
$c.ClearFormats()

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0757'
$c.ClearFormats()

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +1.85%  '
$c.ClearFormats()

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0999'
$c.ClearFormats()

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +1.80%  '
$c.ClearFormats()

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '15.33'
$c.ClearFormats()

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +11.22%  '
$c.ClearFormats()

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.816'
$c.ClearFormats()

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +6.19%  '
$c.ClearFormats()

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +1.73%  '
$c.ClearFormats()

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.10'
$c.ClearFormats()

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +1.92%  '
$c.ClearFormats()

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '1.913.02'
$c.ClearFormats()

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +1.53%  '
$c.ClearFormats()

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '36.956.07'
$c.ClearFormats()

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +4.64%  '
$c.ClearFormats()

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '74.30'
$c.ClearFormats()

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +1.32%  '
$c.ClearFormats()

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.0₃0852'
$c.ClearFormats()

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +3.31%  '
$c.ClearFormats()

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '13.55'
$c.ClearFormats()

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +5.96%  '
$c.ClearFormats()

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '249.75'
$c.ClearFormats()

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +2.40%  '
$c.ClearFormats()

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '5.13'
$c.ClearFormats()

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.ClearFormats()

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.ClearFormats()

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -6.29%  '
$c.ClearFormats()

$c = $ws.Range('B26')
$c.NumberFormat = '@'
$c.Value = 'Monero'
$c.ClearFormats()

$c = $ws.Range('C26')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c.ClearFormats()

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '166.76'
$c.ClearFormats()

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +0.79%  '
$c.ClearFormats()

$c = $ws.Range('B27')
$c.NumberFormat = '@'
$c.Value = 'PancakeSwap'
$c.ClearFormats()

$c = $ws.Range('C27')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c.ClearFormats()

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '2.09'
$c.ClearFormats()

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -3.03%  '
$c.ClearFormats()

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '8.77'
$c.ClearFormats()

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +2.01%  '
$c.ClearFormats()

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '18.60'
$c.ClearFormats()

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +1.93%  '
$c.ClearFormats()

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +0.45%  '
$c.ClearFormats()

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.53'
$c.ClearFormats()

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +5.90%  '
$c.ClearFormats()

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.0607'
$c.ClearFormats()

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +1.82%  '
$c.ClearFormats()

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0906'
$c.ClearFormats()

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +26.79%  '
$c.ClearFormats()

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.26'
$c.ClearFormats()

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +2.53%  '
$c.ClearFormats()

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.91'
$c.ClearFormats()

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +0.93%  '
$c.ClearFormats()

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.16%  '
$c.ClearFormats()

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '18.88'
$c.ClearFormats()

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +37.81%  '
$c.ClearFormats()

$c = $ws.Range('B38')
$c.NumberFormat = '@'
$c.Value = 'ImmutableX'
$c.ClearFormats()

$c = $ws.Range('C38')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.ClearFormats()

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.877'
$c.ClearFormats()

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +3.34%  '
$c.ClearFormats()

$c = $ws.Range('B39')
$c.NumberFormat = '@'
$c.Value = 'TrustWalletToken'
$c.ClearFormats()

$c = $ws.Range('C39')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c.ClearFormats()

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.44'
$c.ClearFormats()

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -1.76%  '
$c.ClearFormats()

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.95'
$c.ClearFormats()

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +1.11%  '
$c.ClearFormats()

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '104.13'
$c.ClearFormats()

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +7.06%  '
$c.ClearFormats()

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0225'
$c.ClearFormats()

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +2.42%  '
$c.ClearFormats()

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '17.40'
$c.ClearFormats()

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +1.79%  '
$c.ClearFormats()

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +19.16%  '
$c.ClearFormats()

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +1.77%  '
$c.ClearFormats()

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.345.68'
$c.ClearFormats()

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +2.43%  '
$c.ClearFormats()

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.ClearFormats()

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +0.30%  '
$c.ClearFormats()

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0834'
$c.ClearFormats()

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +2.88%  '
$c.ClearFormats()

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.80'
$c.ClearFormats()

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +2.24%  '
$c.ClearFormats()

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '6.36'
$c.ClearFormats()

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +1.57%  '
$c.ClearFormats()

$c = $ws.Range('B51')
$c.NumberFormat = '@'
$c.Value = 'RocketPoolETH'
$c.ClearFormats()

$c = $ws.Range('C51')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c.ClearFormats()

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.100.15'
$c.ClearFormats()

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +2.10%  '
$c.ClearFormats()

